$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (YZA567 / 11/10/2024 / Mantenimiento / "" / 12 / En viaje / "" / 15:11:07).
# This shifts the old row 3 (VWX234 ...) up to row 2, carrying its empty-string
# cells (D/G) along unmodified so they stay literal empty-string cells.
$ws.Rows(2).Delete()

# New last row (row 3): DEF567 truck entry. D3/G3 (Descripcion / Tiempo de
# Reparacion) are left blank, same as the other rows' blank entries.
$ws.Range("A3").Value = "DEF567"
$ws.Range("B3").Value = "25/10/2024"
$ws.Range("C3").Value = "Combustible"
$ws.Range("E3").Value = 144
$ws.Range("F3").Value = "FUNCIONAL"
$ws.Range("H3").Value = "09:14:59"
